# "Generate Report for Handback"
#
# 1. The status text shown throughout the workbook flips from
#    "Ready for handoff" to "Handed back: in sync with en-US" (every
#    cell that shows that status, on every sheet).
# 2. Each language sheet (zh-cn / de-de) grows two new populated columns,
#    "Latest Target File" (F) and "Latest Handback File" (G), for both
#    data rows - each gets the a.md source file name and the language's
#    already-handed-off xlf file name, hyperlinked just like the
#    existing columns.
# 3. "Latest Handback DateTime" (H) - previously the placeholder
#    "0001-01-01 00:00:00" - now carries a real timestamp, and the two
#    language sheets disagree (handback finished at a different moment
#    per language).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$zhXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$aMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/35c6247700b445d4ea3063bf1288d20a8da01c8c/e2e/a.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f87f9052c6a330368c97d4bc718f7c83d08eeccc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fcfb83669586c8e95ec8dbba815d8619bae4332e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# --- 1. Flip every "Ready for handoff" cell to the handed-back status ---
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# --- 2. Populate "Latest Target File" (F) / "Latest Handback File" (G) ---
$zh.Hyperlinks.Add($zh.Range("F2"), $aMdUrl, "", "", "a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G2"), $zhXlfUrl, "", "", $zhXlf) | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), $aMdUrl, "", "", "a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G3"), $zhXlfUrl, "", "", $zhXlf) | Out-Null

$de.Hyperlinks.Add($de.Range("F2"), $aMdUrl, "", "", "a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G2"), $deXlfUrl, "", "", $deXlf) | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), $aMdUrl, "", "", "a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G3"), $deXlfUrl, "", "", $deXlf) | Out-Null

# --- 3. Stamp the real handback datetimes (replacing the 0001-01-01 placeholder) ---
$zh.Range("H2").Value = "2016-03-22 12:38:11"
$zh.Range("H3").Value = "2016-03-22 12:38:11"

$de.Range("H2").Value = "2016-03-22 12:38:19"
$de.Range("H3").Value = "2016-03-22 12:38:19"
